# Update cryptos list — refresh prices and 1h volume percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.764.78"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "2.076.85"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.02"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  +3.20%  "
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("D12").Value = "2.375.54"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.96"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.26"
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "2.078.85"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "37.653.37"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("E19").Value = "  -4.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.58"
$ws.Range("E20").Value = "  +1.68%  "
$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.95"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.15"
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.138"
$ws.Range("E27").Value = "  +9.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.91"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.39"
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.64"
$ws.Range("E32").Value = "  +3.85%  "
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.50"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("E37").Value = "  +5.42%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  -4.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0998"
$ws.Range("E40").Value = "  +7.56%  "
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.48"
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("D44").Value = "1.449.01"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.15"
$ws.Range("E46").Value = "  -3.93%  "
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.59"
$ws.Range("E48").Value = "  +2.89%  "
$ws.Range("E49").Value = "  +3.62%  "
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.92"
$ws.Range("E51").Value = "  +7.16%  "
